# Fixed some mistakes in the ppt
#  1. Remove the duplicate "RC5 Timing" slide (originally slide #2 - an
#     exact duplicate of the final "RC5 Timing" slide).
#  2. Merge the two text runs describing the "btn0 not pressed" bullet
#     into a single run (same combined text, same formatting) on the
#     "Buttons & Switches" slide.
#  3. Nudge a few shapes on the "Buttons & Switches" (Input Ukey) slide
#     slightly to fix their alignment.

$p = $ppt.ActivePresentation

# 1) Delete the duplicate "RC5 Timing" slide (originally slide 2).
$p.Slides.Item(2).Delete()

# 2) Merge "If btn0 " + "is not pressed, ..." into one run. That slide
#    ("Buttons & Switches", with the btn0 pseudo-code bullets) is now
#    slide 4 after the deletion above.
$s4 = $p.Slides.Item(4)
$fixedShape = $null
foreach ($shp in $s4.Shapes) {
    if ($shp.Name -eq "Text Box 5") {
        $fixedShape = $shp
        break
    }
}

$tr = $fixedShape.TextFrame.TextRange
$bulletPara = $tr.Paragraphs(1, 1)
$mergedText = "If btn0 is not pressed, r31 = '0'. BNE becomes True. The code will keep looping."
# A same-value assignment is treated as a no-op and would leave the
# original two runs split in the underlying XML, so force a real
# content change first, then set the final (merged) text so it
# collapses into a single run with the original formatting.
$bulletPara.Text = "__tmp__"
$bulletPara = $tr.Paragraphs(1, 1)
$bulletPara.Text = $mergedText

# 3) Reposition shapes on the "Buttons & Switches" (Input Ukey) slide,
#    which is now slide 5 after the deletion above.
$s5 = $p.Slides.Item(5)

foreach ($shp in $s5.Shapes) {
    if ($shp.Name -eq "Picture 69") {
        $shp.Left = 101.5
    }
    elseif ($shp.Name -eq "Text Box 4") {
        $shp.Left = 98.75
        $shp.Top = 299.90001
    }
    elseif ($shp.Name -eq "Picture 7") {
        $shp.Left = 148.15001
    }
}
